$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value2 = "27.581.11"
$ws.Cells.Item(2,5).Value2 = "  -2.46%  "

# Row 3
$ws.Cells.Item(3,4).Value2 = "1.840.70"
$ws.Cells.Item(3,5).Value2 = "  -1.52%  "

# Row 4
$ws.Cells.Item(4,5).Value2 = "  -0.31%  "

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = "313.56"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value2 = "  -1.69%  "

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = "1.000"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value2 = "  -0.27%  "

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = "0.4220"
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).Value2 = "  -4.06%  "

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value2 = "0.3627"
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value2 = "  -1.68%  "

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = "45.23"
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value2 = "  +0.21%  "

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = "0.07240"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value2 = "  -3.25%  "

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = "0.8872"
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).Value2 = "  -5.18%  "

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = "20.55"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value2 = "  -3.88%  "

# Row 13
$ws.Cells.Item(13,4).Value2 = "1.842.81"
$ws.Cells.Item(13,5).Value2 = "  -2.86%  "

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = "5.360"
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).Value2 = "  -1.81%  "

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = "6.541"
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).Value2 = "  -2.24%  "

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value2 = "0.06853"
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value2 = "  -0.74%  "

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value2 = "1.002"
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value2 = "  -0.21%  "

# Row 18
$ws.Cells.Item(18,5).Value2 = "  -4.39%  "

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value2 = "0.000008793"
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value2 = "  -2.52%  "

# Row 20
$ws.Cells.Item(20,5).Value2 = "  -0.28%  "

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = "15.43"
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value2 = "  -3.13%  "

# Row 22
$ws.Cells.Item(22,4).Value2 = "27.567.69"
$ws.Cells.Item(22,5).Value2 = "  -2.48%  "

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = "4.982"
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value2 = "  -2.56%  "

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = "10.51"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value2 = "  -2.78%  "

# Row 25
$ws.Cells.Item(25,4).Value2 = "2.056.93"
$ws.Cells.Item(25,5).Value2 = "  -2.84%  "

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value2 = "2.032"
$ws.Cells.Item(26,4).ClearFormats()

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = "154.66"
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).Value2 = "  -0.38%  "

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value2 = "18.41"
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value2 = "  +0.50%  "

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = "118.56"
$ws.Cells.Item(29,4).ClearFormats()
$ws.Cells.Item(29,5).Value2 = "  +4.64%  "

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = "5.197"
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value2 = "  -2.21%  "

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value2 = "1.802"
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value2 = "  +4.92%  "

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value2 = "0.08862"
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value2 = "  -1.98%  "

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = "0.7716"
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value2 = "  -2.86%  "

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = "4.551"
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value2 = "  -5.97%  "

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = "2.951"
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value2 = "  +0.90%  "

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value2 = "1.099"
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value2 = "  -6.18%  "

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = "0.9996"
$ws.Cells.Item(37,4).ClearFormats()
$ws.Cells.Item(37,5).Value2 = "  -0.36%  "

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value2 = "0.05398"
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value2 = "  -0.81%  "

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = "1.095"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value2 = "  -2.84%  "

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = "0.01912"
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value2 = "  -2.93%  "

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = "2.780"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value2 = "  -6.71%  "

# Row 42
$ws.Cells.Item(42,2).Value2 = "FraxShare"
$ws.Cells.Item(42,3).Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = "6.804"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value2 = "  -3.95%  "

# Row 43
$ws.Cells.Item(43,2).Value2 = "TheSandbox"
$ws.Cells.Item(43,3).Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value2 = "0.5039"
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value2 = "  -4.21%  "

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value2 = "0.1643"
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value2 = "  -2.29%  "

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = "0.06607"
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).Value2 = "  -2.13%  "

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = "8.138"
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value2 = "  -6.58%  "

# Row 47
$ws.Cells.Item(47,2).Value2 = "Decentraland"
$ws.Cells.Item(47,3).Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = "0.4679"
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value2 = "  -3.89%  "

# Row 48
$ws.Cells.Item(48,2).Value2 = "Quant"
$ws.Cells.Item(48,3).Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value2 = "104.96"
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).Value2 = "  -2.03%  "

# Row 49
$ws.Cells.Item(49,2).Value2 = "EnergySwap"
$ws.Cells.Item(49,3).Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = "10.26"
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value2 = "  -2.24%  "

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = "0.9997"
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value2 = "  -0.32%  "

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = "1.620"
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).Value2 = "  -3.06%  "

